# Add "NA" values under the duplicate_image_filename column (column E)
# for the data rows (rows 2 through 21) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
